$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 100m - F (Final) results row (row 3)
$ws.Range("B3").Value = "Julien Alfred"
$ws.Range("C3").Value = "Sha'Carri Richardson"
$ws.Range("D3").Value = "Melissa Jefferson"

# Move the active selection to D9 (matches the commit's final cursor position)
$ws.Range("D9").Select()
